{"js": "// Fix navigation between Asset windows.\n//\n// 1) Give the \"{% for annotation in annotations %}\" paragraph a hanging\n//    indent (left 720 twips / 36pt, hanging 720 twips / 36pt) so wrapped\n//    continuation lines line up under the loop body instead of the tag.\n// 2) Extend the nested \"{% for schange in annotation.statechanges %}\"\n//    Jinja tag with an \"if loop.first\" clause so only the first state\n//    change is iterated (the actual navigation fix).\n\n// --- 1) Hanging indent on the outer \"{% for annotation in annotations %}\" paragraph ---\nconst forAnnotation = context.document.body.search(\n  \"{% for annotation in annotations %}\",\n  { matchCase: true }\n);\nforAnnotation.load(\"paragraphs\");\nawait context.sync();\n\nconst annotationPara = forAnnotation.items[0].paragraphs.getFirst();\n// 36 points == 720 twips (Word stores w:ind in twentieths of a point).\nannotationPara.leftIndent = 36;\nannotationPara.firstLineIndent = -36;\nawait context.sync();\n\n// --- 2) Insert \" if loop.first\" right after \"annotation.statechanges\" ---\nconst stateChanges = context.document.body.search(\n  \"annotation.statechanges\",\n  { matchCase: true }\n);\nstateChanges.load(\"text\");\nawait context.sync();\n\nstateChanges.items[0].insertText(\" if loop.first\", \"After\");\nawait context.sync();\n", "ps1": "# Fix navigation between Asset windows.\n#\n# 1) Give the \"{% for annotation in annotations %}\" paragraph a hanging\n#    indent (left 720 twips / 36pt, hanging 720 twips / 36pt) so wrapped\n#    continuation lines line up under the loop body instead of the tag.\n# 2) Extend the nested \"{% for schange in annotation.statechanges %}\"\n#    Jinja tag with an \"if loop.first\" clause so only the first state\n#    change is iterated (the actual navigation fix).\n\n$d = $word.ActiveDocument\n\n# --- 1) Hanging indent on the outer \"{% for annotation in annotations %}\" paragraph ---\n$range1 = $d.Content\n$null = $range1.Find.Execute(\"{% for annotation in annotations %}\")\n$para = $range1.Paragraphs.First\n$para.Format.LeftIndent = 36\n$para.Format.FirstLineIndent = -36\n\n# --- 2) Insert \" if loop.first\" right after \"annotation.statechanges\" ---\n$range2 = $d.Content\n$null = $range2.Find.Execute(\"annotation.statechanges\")\n$range2.Collapse(0)\n$range2.InsertAfter(\" if loop.first\")\n"}
